$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "Periodo Mora" values for the 3 existing workers (2507->2509)
#    and recompute the dependent summary cells.
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "2509"
$ws.Range("E17").Value = "2509"
$ws.Range("E18").Value = "2509"

# "VALOR MORA" total (B11 label / E11 value) goes from 341640 to 170820
$ws.Range("E11").Value = 170820

# "Cant. Periodos" goes from 2 to 1 (one period remains: 2509)
$ws.Range("F13").Value = 1

# ---------------------------------------------------------------------------
# 2) Remove the old second-period rows (19:21) that held the "2508" entries -
#    the account statement now only lists a single period per worker.
# ---------------------------------------------------------------------------
$ws.Rows("19:21").Delete()

# ---------------------------------------------------------------------------
# 3) Center the "Periodo Mora" column values (E16:E18) - matches the new
#    layout applied when the table was rebuilt for the remaining rows.
# ---------------------------------------------------------------------------
$ws.Range("E16:E18").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4) Give the new last data row (18) the closing bottom border of the table,
#    same treatment the previous last row (21) used to have.
# ---------------------------------------------------------------------------
$closing = $ws.Range("B18:J18")
$closing.Borders.Item(9).LineStyle = 1
$closing.Borders.Item(9).Weight = 2
$closing.Borders.Item(9).ColorIndex = 0

Write-Output "done"
